# Insert a new row at the top of the data table, shifting all existing
# rows (1-31) down by one (new rows 2-32), and populate the new row 1
# with the values 0, 0, 0, 40 in columns A-D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(1).Insert()

$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 0
$ws.Range("D1").Value = 40

$ws.Range("A1:XFD1").Select()
